$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Tocantins"
$ws.Range("B2").Value = "Diferença 2025/04 - 2024/04"
$ws.Range("C2").Value = 2.17

# Row 3
$ws.Range("A3").Value = "Amazonas"
$ws.Range("B3").Value = "Diferença 2025/04 - 2024/04"
$ws.Range("C3").Value = 2.05

# Row 4
$ws.Range("A4").Value = "Bahia"
$ws.Range("B4").Value = "Diferença 2025/04 - 2024/04"
$ws.Range("C4").Value = 1.98

# Row 5
$ws.Range("A5").Value = "Acre"
$ws.Range("B5").Value = "Diferença 2025/04 - 2024/04"
$ws.Range("C5").Value = 1.95

# Row 6
$ws.Range("B6").Value = "Diferença 2025/04 - 2024/04"
$ws.Range("C6").Value = 0.84

# Row 7
$ws.Range("A7").Value = "Maranhão"
$ws.Range("B7").Value = "Diferença 2025/04 - 2024/04"
$ws.Range("C7").Value = 0.78

# Row 8
$ws.Range("B8").Value = "Diferença 2025/04 - 2024/04"
$ws.Range("C8").Value = -2.01
$ws.Range("D8").Value = "26º"

# Row 9
$ws.Range("B9").Value = "Diferença 2025/04 - 2024/04"
$ws.Range("C9").Value = 0.41

# Row 10
$ws.Range("B10").Value = "Diferença 2025/04 - 2024/04"
$ws.Range("C10").Value = 0.73
